$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.158404049600987
$ws.Range("C2").Value = 0.001492798534567874
$ws.Range("B3").Value = 0.3448363805299818
$ws.Range("C3").Value = 0.003933342534632842
